$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Project" sheet: drop the old Sales/Product-analytics bullet lists,
#    keep the Schema/Metadata/Sales Analytics/Product Analytics headers.
# ---------------------------------------------------------------------------
$project = $wb.Worksheets.Item("Project")
$project.Range("C6:C13").ClearContents()
$project.Range("C16:C20").ClearContents()
$project.Range("A6:XFD6").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Add the new "Flow Chart" sheet after "Project" (becomes the active tab).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$flow = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$flow.Name = "Flow Chart"

$flow.Range("B2").Value = "Retailmart Database"

$flow.Range("C3").Value = "8 Schemas"

$flow.Range("C4").Value = "analytics"
$flow.Range("D4").Value = "customers"

$flow.Range("D5").Value = "sales"
$flow.Range("D6").Value = "stores"
$flow.Range("D7").Value = "products"

$flow.Range("F4").Value = "kpi"
$flow.Range("G4").Value = "trends"
$flow.Range("H4").Value = "monthly / weekly / yearly"

$flow.Range("D10").Value = "Customer_anlytics.sql"

$flow.Range("E11").Value = "Views"
$flow.Range("E12").Value = "Materialized Views"
$flow.Range("E13").Value = "refresh_procedure"
$flow.Range("E14").Value = "function convert our data into json files"

$flow.Range("H15").Value = "Front End"
$flow.Range("H17").Value = "style.css"

$flow.Range("E16").Value = "Execute these Functions"
$flow.Range("E17").Value = ".json "

$flow.Range("H16").Value = "index.html"
$flow.Range("H18").Value = "dashboard.js"

# Column widths (best-fit, matching the authored layout). The COM layer
# quantizes ColumnWidth to 1/6-character steps and adds a constant offset
# (5/6 char) between the value assigned and the value persisted in the XML,
# so we pre-compensate to land on the nearest achievable width.
$flow.Columns.Item(2).ColumnWidth = 18.1640625 - (5/6)
$flow.Columns.Item(4).ColumnWidth = 19.33203125 - (5/6)
$flow.Columns.Item(5).ColumnWidth = 32.83203125 - (5/6)
$flow.Columns.Item(6).ColumnWidth = 3.5 - (5/6)
$flow.Columns.Item(7).ColumnWidth = 6.33203125 - (5/6)
$flow.Columns.Item(8).ColumnWidth = 20.83203125 - (5/6)

$flow.Range("F12").Select() | Out-Null
$excel.ActiveWindow.Zoom = 227
